$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 194, pushing the existing rows 194-204 down to 196-206.
$ws.Rows("194:195").Insert()

# New row 194: weekly price entry (Provincia del Elquí, dated 2021-11-16)
$ws.Range("A194").Value = 10
$ws.Range("B194").Value = 'Vega Modelo de Temuco'
$ws.Range("C194").Value = 'La Araucanía'
$ws.Range("D194").Value = 44516
$ws.Range("E194").Value = 9
$ws.Range("F194").Value = 100112017
$ws.Range("G194").Value = 'Apio'
$ws.Range("H194").Value = 'Americana (o)'
$ws.Range("I194").Value = 'Primera'
$ws.Range("J194").Value = 95
$ws.Range("K194").Value = 8000
$ws.Range("L194").Value = 8000
$ws.Range("M194").Value = 8000
$ws.Range("N194").Value = '$/docena de matas'
$ws.Range("O194").Value = 'Provincia del Elquí'
$ws.Range("P194").Value = 1333
$ws.Range("Q194").Value = 6
$ws.Range("R194").Value = 'Hortaliza'

# New row 195: weekly price entry (Región Metropolitana, dated 2021-11-16)
$ws.Range("A195").Value = 10
$ws.Range("B195").Value = 'Vega Modelo de Temuco'
$ws.Range("C195").Value = 'La Araucanía'
$ws.Range("D195").Value = 44516
$ws.Range("E195").Value = 9
$ws.Range("F195").Value = 100112017
$ws.Range("G195").Value = 'Apio'
$ws.Range("H195").Value = 'Americana (o)'
$ws.Range("I195").Value = 'Primera'
$ws.Range("J195").Value = 55
$ws.Range("K195").Value = 8000
$ws.Range("L195").Value = 8000
$ws.Range("M195").Value = 8000
$ws.Range("N195").Value = '$/docena de matas'
$ws.Range("O195").Value = 'Región Metropolitana'
$ws.Range("P195").Value = 1333
$ws.Range("Q195").Value = 6
$ws.Range("R195").Value = 'Hortaliza'
